# Move the last 35 comment rows from Sheet1 (A271:A305) up to the top of
# Sheet2 (A1:A35), pushing Sheet2's existing rows down. This mirrors
# "taking comments between workbook sheets": the comments disappear from
# the bottom of Sheet1 and reappear at the top of Sheet2.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Make room at the top of Sheet2 for the 35 rows coming from Sheet1.
$ws2.Range("A1:A35").Insert()

# Copy the trailing comments out of Sheet1 into the newly freed space.
$ws1.Range("A271:A305").Copy()
$ws2.Range("A1:A35").PasteSpecial()

# Remove the now-duplicated rows from the bottom of Sheet1.
$ws1.Range("A271:A305").EntireRow.Delete()
